$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Crypto price/volume refresh (GitHub Actions symbol-list update).
# Columns: D = Price, E = Volume(1h). Values are stored as plain text
# (as in the source data), so NumberFormat is forced to "@" (Text)
# before the write to stop Excel from auto-coercing numeric-looking
# strings (and percentages) into Number/Percentage cells, then the
# temporary format override is cleared so no stray style lingers.
function Set-TextValue($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

# Row 2
Set-TextValue 2 4 "300.73"
Set-TextValue 2 5 "0.32%"
# Row 3
Set-TextValue 3 4 "31.48"
Set-TextValue 3 5 "-0.07%"
# Row 4
Set-TextValue 4 4 "5.081"
Set-TextValue 4 5 "-1.04%"
# Row 5
Set-TextValue 5 4 "0.07862"
Set-TextValue 5 5 "-2.01%"
# Row 6
Set-TextValue 6 4 "2.330"
Set-TextValue 6 5 "-3.00%"
# Row 7
Set-TextValue 7 4 "7.808"
Set-TextValue 7 5 "-1.29%"
# Row 8
Set-TextValue 8 4 "3.830"
Set-TextValue 8 5 "-0.50%"
# Row 9
Set-TextValue 9 4 "0.9208"
Set-TextValue 9 5 "0.82%"
# Row 10
Set-TextValue 10 4 "0.1754"
Set-TextValue 10 5 "1.43%"
# Row 11
Set-TextValue 11 4 "0.07577"
Set-TextValue 11 5 "3.09%"
# Row 12
Set-TextValue 12 4 "0.09193"
Set-TextValue 12 5 "14.43%"
# Row 13
Set-TextValue 13 4 "0.03001"
Set-TextValue 13 5 "-1.24%"
# Row 14
Set-TextValue 14 4 "0.1003"
Set-TextValue 14 5 "0.58%"
# Row 15
Set-TextValue 15 4 "0.001504"
Set-TextValue 15 5 "-0.54%"
# Row 16
Set-TextValue 16 4 "0.005824"
Set-TextValue 16 5 "-4.92%"
# Row 17
Set-TextValue 17 5 "-0.70%"
# Row 18
Set-TextValue 18 5 "-0.83%"
# Row 19
Set-TextValue 19 5 "-0.93%"
# Row 20
Set-TextValue 20 4 "0.1296"
Set-TextValue 20 5 "-3.96%"
# Row 21
Set-TextValue 21 5 "-11.55%"
# Row 22
Set-TextValue 22 4 "0.1709"
Set-TextValue 22 5 "6.27%"
# Row 23
Set-TextValue 23 4 "0.04596"
Set-TextValue 23 5 "-0.41%"
# Row 24
Set-TextValue 24 5 "-1.11%"
# Row 25
Set-TextValue 25 4 "0.004469"
Set-TextValue 25 5 "0.40%"
# Row 26
Set-TextValue 26 4 "0.0001249"
Set-TextValue 26 5 "5.16%"
# Row 27
Set-TextValue 27 4 "0.0003382"
Set-TextValue 27 5 "-1.98%"
# Row 39
Set-TextValue 39 4 "0.01748"
Set-TextValue 39 5 "-4.39%"
# Row 40
Set-TextValue 40 4 "0.04731"
Set-TextValue 40 5 "4.32%"
# Row 41
Set-TextValue 41 4 "0.007045"
Set-TextValue 41 5 "-3.36%"
# Row 42
Set-TextValue 42 4 "0.1360"
Set-TextValue 42 5 "1.03%"
# Row 43
Set-TextValue 43 4 "0.002188"
Set-TextValue 43 5 "0.05%"
# Row 44
Set-TextValue 44 4 "0.009759"
Set-TextValue 44 5 "-8.32%"
# Row 45
Set-TextValue 45 4 "0.00006277"
Set-TextValue 45 5 "-0.19%"
# Row 46
Set-TextValue 46 4 "0.00000000749"
Set-TextValue 46 5 "-0.72%"
# Row 47
Set-TextValue 47 5 "19.50%"
# Row 48
Set-TextValue 48 4 "1.153"
Set-TextValue 48 5 "40.48%"
# Row 49
Set-TextValue 49 4 "0.00002098"
Set-TextValue 49 5 "-0.72%"
# Row 50
Set-TextValue 50 4 "0.0001998"
Set-TextValue 50 5 "-0.72%"
